$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the diff. D and E columns hold numeric-looking
# text (prices / percentages) that must remain plain text, so we force the
# "@" text number format before assigning, then clear the format again so the
# cell keeps its original (unstyled) appearance.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "61.699.19"
Set-TextValue "E2" "  -5.36%  "
Set-TextValue "D3" "3.118.00"
Set-TextValue "E3" "  -7.68%  "
Set-TextValue "E4" "  -0.30%  "
Set-TextValue "D5" "509.95"
Set-TextValue "E5" "  -3.12%  "
Set-TextValue "D6" "166.18"
Set-TextValue "E6" "  -9.39%  "
Set-TextValue "D7" "0.580"
Set-TextValue "E7" "  -3.25%  "
Set-TextValue "E8" "  -0.02%  "
Set-TextValue "D9" "3.116.97"
Set-TextValue "E9" "  -7.55%  "
Set-TextValue "D10" "0.583"
Set-TextValue "E10" "  -6.31%  "
Set-TextValue "D11" "51.00"
Set-TextValue "E11" "  -10.57%  "
Set-TextValue "D12" "0.125"
Set-TextValue "E12" "  -5.25%  "
Set-TextValue "D13" "0.0000242"
Set-TextValue "E13" "  -3.93%  "
Set-TextValue "D14" "8.71"
Set-TextValue "E14" "  -5.50%  "
Set-TextValue "D15" "3.603.04"
Set-TextValue "E15" "  -8.05%  "
Set-TextValue "D16" "0.113"
Set-TextValue "E16" "  -7.34%  "
Set-TextValue "D17" "3.108.89"
Set-TextValue "D18" "61.361.55"
Set-TextValue "E18" "  -5.56%  "
Set-TextValue "D19" "16.60"
Set-TextValue "E19" "  -3.98%  "
Set-TextValue "D20" "10.62"
Set-TextValue "E20" "  -3.45%  "
Set-TextValue "D21" "0.937"
Set-TextValue "E21" "  -2.46%  "
Set-TextValue "D22" "354.93"
Set-TextValue "E22" "  -4.26%  "
$ws.Range("B23").Value = "RenderToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D23" "11.02"
Set-TextValue "E23" "  +2.31%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D24" "78.71"
Set-TextValue "E24" "  -2.75%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D25" "3.60"
Set-TextValue "E25" "  -2.45%  "
Set-TextValue "E26" "  +4.46%  "
Set-TextValue "D27" "3.77"
Set-TextValue "E27" "  +0.85%  "
Set-TextValue "D28" "2.53"
Set-TextValue "E28" "  -3.50%  "
Set-TextValue "D29" "10.82"
Set-TextValue "E29" "  -4.82%  "
Set-TextValue "D30" "7.86"
Set-TextValue "E30" "  -6.92%  "
Set-TextValue "D31" "632.58"
Set-TextValue "E31" "  -5.36%  "
Set-TextValue "D32" "27.41"
Set-TextValue "E32" "  -7.00%  "
Set-TextValue "D33" "6.18"
Set-TextValue "E33" "  -7.76%  "
Set-TextValue "D34" "10.97"
Set-TextValue "E34" "  -0.92%  "
Set-TextValue "E35" "  -0.01%  "
Set-TextValue "D36" "0.101"
Set-TextValue "E36" "  -3.00%  "
Set-TextValue "D37" "55.66"
Set-TextValue "E37" "  -8.93%  "
Set-TextValue "D38" "34.87"
Set-TextValue "E38" "  -3.66%  "
Set-TextValue "D39" "0.364"
Set-TextValue "E39" "  -3.51%  "
Set-TextValue "D40" "0.998"
Set-TextValue "E40" "  -0.13%  "
Set-TextValue "D41" "0.0₃0661"
Set-TextValue "E41" "  +5.76%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D42" "2.50"
Set-TextValue "E42" "  +7.00%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D43" "0.119"
Set-TextValue "E43" "  -6.46%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D44" "2.792.40"
Set-TextValue "E44" "  -1.04%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D45" "2.85"
Set-TextValue "E45" "  +11.09%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D46" "2.60"
Set-TextValue "E46" "  -0.79%  "
Set-TextValue "D47" "0.0377"
Set-TextValue "E47" "  -2.76%  "
Set-TextValue "D48" "2.93"
Set-TextValue "E48" "  +1.76%  "
Set-TextValue "D49" "2.44"
Set-TextValue "E49" "  -11.40%  "
Set-TextValue "D50" "0.120"
Set-TextValue "E50" "  -3.17%  "
Set-TextValue "D51" "129.21"
Set-TextValue "E51" "  -4.77%  "
